# Limpieza de separadores y parrafos vacios.
#
# Elimina (en todo el documento):
#   1. Los parrafos que contienen imagenes en linea (InlineShapes).
#   2. Los parrafos "separador" formados solo por una linea de guiones
#      Unicode (U+2500 BOX DRAWINGS LIGHT HORIZONTAL).
#   3. Los parrafos vacios con "spacing before = 40 twips" (2 pt) que
#      aparecian justo despues de cada tabla.

$d = $word.ActiveDocument

# Recolectar los INDICES (no las posiciones de caracter) de los
# parrafos de $d.Paragraphs que hay que eliminar. Se recolectan antes
# de borrar nada porque, al eliminar un parrafo, los indices de los
# parrafos anteriores a el en la coleccion no cambian -- solo cambian
# (se reducen) los indices de los parrafos posteriores. Por eso el
# borrado se hace despues, recorriendo los indices de mayor a menor.
$indices = New-Object System.Collections.ArrayList

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $len = $rng.End - $rng.Start
    $delete = $false

    if ($len -gt 0) {
        $txt = $rng.Text
        $firstChar = $txt.Substring(0, 1)

        if ($firstChar -eq [char]0x2500) {
            # Parrafo separador "────────────────────"
            $delete = $true
        } elseif ($rng.Information(12) -eq $false -and $p.Format.SpaceBefore -eq 2 -and $len -eq 1) {
            # Parrafo vacio (spacing before = 40 twips = 2 pt), fuera de
            # cualquier tabla, formado solo por la marca de parrafo.
            $delete = $true
        }
    }

    if (-not $delete -and $rng.InlineShapes.Count -gt 0) {
        # Parrafo que solo contiene una imagen en linea
        $delete = $true
    }

    if ($delete) {
        [void]$indices.Add($i)
    }
}

# Borrar de mayor indice a menor para que los indices ya recolectados
# sigan apuntando a los parrafos correctos.
$sortedIdx = $indices | Sort-Object -Descending
foreach ($idx in $sortedIdx) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

Write-Output "Parrafos eliminados: $($sortedIdx.Count)"
Write-Output "Parrafos restantes: $($d.Paragraphs.Count)"
